$d = $word.ActiveDocument

$replacements = @(
    @("58÷4=", "89÷3="),
    @("94÷3=", "96÷2="),
    @("87÷8=", "32÷8="),
    @("35÷3=", "18÷7="),
    @("96÷9=", "49÷5="),
    @("20÷6=", "43÷3="),
    @("99÷6=", "84÷6="),
    @("70÷8=", "69÷4="),
    @("74÷5=", "57÷2="),
    @("94÷6=", "98÷3="),
    @("23÷7=", "88÷3="),
    @("34÷7=", "27÷6="),
    @("91÷3=", "87÷4="),
    @("71÷8=", "98÷3="),
    @("22÷6=", "94÷2="),
    @("84÷9=", "89÷3="),
    @("38÷5=", "76÷4="),
    @("25÷6=", "48÷3="),
    @("23÷5=", "44÷6="),
    @("74÷2=", "10÷9="),
    @("43÷4=", "39÷6="),
    @("67÷4=", "98÷9="),
    @("80÷2=", "46÷6="),
    @("36÷9=", "51÷7="),
    @("53÷3=", "45÷7=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}
